# Add "Wins" / "Losses" / "Ties" season-record columns (AC, AD, AE)
# to the player table on Sheet1, matching the header style used by the
# existing header row and filling every data row (2-42) with the
# team's 1996 season record (90 wins, 72 losses, 0 ties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the existing header formatting (bold, border, centered) from A1
# onto the new header cells so they match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# --- Data rows -----------------------------------------------------------
$wins = 90
$losses = 72
$ties = 0

for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 29).Value = $wins
    $ws.Cells.Item($row, 30).Value = $losses
    $ws.Cells.Item($row, 31).Value = $ties
}
